$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns refreshed by the scheduled GitHub Actions run.
# D values are forced to Text (NumberFormat "@") before assignment, then the style is
# reset back to "Normal" so only the cell value changes (matches the source data which
# stores these as plain text, not numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.088.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.912.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.80%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.70%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4829'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3821'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9362'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07826'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.913.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.497'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.629'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.20'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.009'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008829'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.121.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.158'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.154.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.097'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.958'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08911'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.370'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.248'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7685'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.684'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.625'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02045'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.098'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05314'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5504'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.995'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.040'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1524'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.440'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4846'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.656'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.68%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.40'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06101'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '

# Row 26 only has its Volume(1h) figure updated; Price stays 1.920
$ws.Range("E26").Value = '  -2.36%  '

